$d = $word.ActiveDocument

function Get-CellRange($cell) {
    # Build a fresh Document.Range bound to the cell's character extents.
    # (Using the Cell's own .Range object directly with .Find can cause the
    # search to run against the whole document instead of being confined to
    # the cell, so we re-derive an explicit Range from its Start/End.)
    $inner = $cell.Range
    return $d.Range($inner.Start, $inner.End)
}

# --- Change 1: Table 1 (header info table), Row 2, Col 2 ("DIA" value): "    23" -> "    26"
$t1 = $d.Tables.Item(1)
$r1 = Get-CellRange $t1.Cell(2, 2)
$ok1 = $r1.Find.Execute("    23", $true, $false, $false, $false, $false, $true, 0, $false, "    26", 1)
Write-Host "Change 1 (DIA 23 -> 26):" $ok1

# --- Change 2: Table 2 (articulos table), Row 1, Col 2: "RECLUTAMIENTO" -> "SISETMASPRUEBA"
$t2 = $d.Tables.Item(2)
$r2 = Get-CellRange $t2.Cell(1, 2)
$ok2 = $r2.Find.Execute("RECLUTAMIENTO", $true, $false, $false, $false, $false, $true, 0, $false, "SISETMASPRUEBA", 1)
Write-Host "Change 2 (RECLUTAMIENTO -> SISETMASPRUEBA):" $ok2

# --- Change 3: Table 2, Row 3, Col 2 (name): "omar" -> "noe"
$r3 = Get-CellRange $t2.Cell(3, 2)
$ok3 = $r3.Find.Execute("omar", $true, $false, $false, $false, $false, $true, 0, $false, "noe", 1)
Write-Host "Change 3 (omar -> noe):" $ok3

# --- Change 4: Table 2, Row 3, Col 4 (CANTIDAD): "1" -> "2"
$r4 = Get-CellRange $t2.Cell(3, 4)
$ok4 = $r4.Find.Execute("1", $true, $false, $false, $false, $false, $true, 0, $false, "2", 1)
Write-Host "Change 4 (CANTIDAD 1 -> 2):" $ok4

$d.Save()
